$d = $word.ActiveDocument

# --- 1. Text edits inside the "Design brief" paragraph -------------------

# Insert "during the school holidays " after "...friends over " and before
# "as long as...".
$d.Content.Find.Execute(
    "some friends over as long as", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "some friends over during the school holidays as long as", 2) | Out-Null

# Insert " at any given time" after "...in the house" and before the period
# that starts "However, ...".
$d.Content.Find.Execute(
    "in the house. However", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "in the house at any given time. However", 2) | Out-Null

# Insert "prototype " before "a simple counter".
$d.Content.Find.Execute(
    "decide to make a simple counter", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "decide to prototype a simple counter", 2) | Out-Null

# Replace "the visitors" with "your friends".
$d.Content.Find.Execute(
    "count the visitors as they come and go", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "count your friends as they come and go", 2) | Out-Null

# --- 2. Remove the stray empty paragraph right after the brief -----------

$d.Paragraphs(3).Range.Delete()

# --- 3. Relocate the "_GoBack" bookmark to sit right after "Pseudocode" --

$headingRange = $d.Content
$found = $headingRange.Find.Execute("Pseudocode", $true, $false, $false,
                                     $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $afterHeading = $d.Range($headingRange.End, $headingRange.End)
    $d.Bookmarks.Add("_GoBack", $afterHeading) | Out-Null
}

# --- 4. Shrink the page margins from 1" to 0.5" (1440 -> 720 twips) ------

$d.PageSetup.TopMargin = 36
$d.PageSetup.BottomMargin = 36
$d.PageSetup.LeftMargin = 36
$d.PageSetup.RightMargin = 36
